$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "14.60", "3.00").
# Force text formatting before assigning so Excel keeps them as literal
# strings (matching the source workbook) instead of silently coercing
# them to numbers and dropping significant trailing zeros.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.228.23'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.129.17'
$ws.Range('E3').Value = '  +0.55%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '604.43'
$ws.Range('E5').Value = '  -0.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.82'
$ws.Range('E6').Value = '  -1.61%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.124.32'
$ws.Range('E8').Value = '  +0.54%  '
$ws.Range('E9').Value = '  +0.47%  '
$ws.Range('E10').Value = '  -0.18%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.39'
$ws.Range('E11').Value = '  +2.95%  '
$ws.Range('E12').Value = '  -0.20%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000255'
$ws.Range('E13').Value = '  +2.22%  '
$ws.Range('E14').Value = '  -0.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.647.57'
$ws.Range('E15').Value = '  +0.92%  '
$ws.Range('E16').Value = '  +2.70%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.173.32'
$ws.Range('E17').Value = '  +0.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.142.60'
$ws.Range('E18').Value = '  +1.14%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.89'
$ws.Range('E19').Value = '  +0.97%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '480.38'
$ws.Range('E20').Value = '  +0.83%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.60'
$ws.Range('E21').Value = '  -0.24%  '
$ws.Range('E22').Value = '  +1.60%  '
$ws.Range('E23').Value = '  -0.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.63'
$ws.Range('E24').Value = '  +2.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.42'
$ws.Range('E25').Value = '  -0.98%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E27').Value = '  -0.84%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.33'
$ws.Range('E28').Value = '  -0.97%  '
$ws.Range('E29').Value = '  +7.05%  '
$ws.Range('E30').Value = '  -3.73%  '
$ws.Range('E31').Value = '  -1.64%  '
$ws.Range('E33').Value = '  +2.47%  '
$ws.Range('E34').Value = '  -2.73%  '
$ws.Range('E35').Value = '  -1.47%  '
$ws.Range('E36').Value = '  +0.73%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0₃0769'
$ws.Range('E37').Value = '  +3.42%  '
$ws.Range('E38').Value = '  -0.87%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.00'
$ws.Range('E39').Value = '  +1.97%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '445.64'
$ws.Range('E40').Value = '  -3.38%  '
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('E42').Value = '  +0.90%  '
$ws.Range('E43').Value = '  -1.68%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.853.98'
$ws.Range('E44').Value = '  +0.42%  '
$ws.Range('E45').Value = '  -2.27%  '
$ws.Range('E46').Value = '  -1.82%  '
$ws.Range('E47').Value = '  +0.75%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.98'
$ws.Range('E49').Value = '  -0.66%  '
$ws.Range('E50').Value = '  +0.25%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.73'
$ws.Range('E51').Value = '  +2.04%  '
